$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -21.371
$ws.Range("A21").Value = -20.88
$ws.Range("A23").Value = -21.368
$ws.Range("A25").Value = -21.937
